$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.250.71'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '3.690.91'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '681.79'
$ws.Range('E5').Value = '  -3.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.48'
$ws.Range('E6').Value = '  -4.46%  '
$ws.Range('D7').Value = '3.689.77'
$ws.Range('E7').Value = '  -2.80%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('E10').Value = '  -7.05%  '
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('E13').Value = '  -6.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.54'
$ws.Range('E14').Value = '  -6.94%  '
$ws.Range('D15').Value = '4.313.47'
$ws.Range('E15').Value = '  -2.83%  '
$ws.Range('D16').Value = '3.691.04'
$ws.Range('E16').Value = '  -2.73%  '
$ws.Range('D17').Value = '69.341.60'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.36'
$ws.Range('E19').Value = '  -5.89%  '
$ws.Range('E20').Value = '  -7.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '484.53'
$ws.Range('E21').Value = '  -2.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.80'
$ws.Range('E22').Value = '  -7.63%  '
$ws.Range('E23').Value = '  -8.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.98'
$ws.Range('E24').Value = '  -5.67%  '
$ws.Range('D25').Value = '3.836.95'
$ws.Range('E25').Value = '  -2.85%  '
$ws.Range('E26').Value = '  -10.94%  '
$ws.Range('E27').Value = '  -4.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('E29').Value = '  -7.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.83'
$ws.Range('E30').Value = '  -10.28%  '
$ws.Range('E31').Value = '  -10.66%  '
$ws.Range('E32').Value = '  -4.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.83'
$ws.Range('E33').Value = '  -6.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.07'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.167'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').Value = '3.653.35'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.55'
$ws.Range('E38').Value = '  -5.78%  '
$ws.Range('E39').Value = '  +1.73%  '
$ws.Range('E40').Value = '  -6.98%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.963'
$ws.Range('E44').Value = '  -7.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '157.48'
$ws.Range('E45').Value = '  -4.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.21'
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.83'
$ws.Range('E47').Value = '  -13.37%  '
$ws.Range('E48').Value = '  -13.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '391.09'
$ws.Range('E49').Value = '  -8.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.09'
$ws.Range('E50').Value = '  -5.83%  '
$ws.Range('E51').Value = '  -4.98%  '
